$wb = $excel.ActiveWorkbook

# --- Update "Bugs and Quirks" sheet ---
# New bug: "Player moves slightly when standing still." goes into C4.
# Existing bug text in B6 is replaced with the new by-pass nextFire bug,
# and the old B7 text is removed (its content moves up conceptually, and
# the remaining text in B7 is no longer present).
$bugs = $wb.Worksheets.Item("Bugs and Quirks")
$bugs.Range("C4").Value = "Player moves slightly when standing still."
$bugs.Range("B6").Value = "Player can by-pass nextFire wait time on sniper rifle by quickly switching to another weapon and switching back and firing."
$bugs.Range("B7").Value = ""
$bugs.Range("B6").Select()

# --- Update "Immediate Checklist" sheet ---
# Remove the old rows 3-5 (completed/obsolete items) and promote the
# sniper rifle by-pass bug (now also duplicated into Bugs&Quirks above)
# into row 2 of the checklist, replacing the old "sniper purchase button" item.
$checklist = $wb.Worksheets.Item("Immediate Checklist")
$checklist.Range("A3:A5").EntireRow.Delete()
$checklist.Range("A2").Value = "Player can by-pass nextFire wait time on sniper rifle by quickly switching to another weapon and switching back and firing."
$checklist.Range("A2").WrapText = $true
$checklist.Rows.Item(2).RowHeight = 30.75

# --- Reorder sheet tabs ---
# Old order: Immediate Checklist, Features, Bugs and Quirks
# New order: Features, Bugs and Quirks, Immediate Checklist
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$checklist.Move($null, $lastSheet)

# Immediate Checklist becomes the active/selected tab. Re-fetch it by name
# since the move can invalidate the earlier object reference.
$wb.Worksheets.Item("Immediate Checklist").Activate()
